# Preisalarm.pptx edits:
#  1) Slide 1  - header date textbox: "17.07.2025" -> "17. Juli 2025"
#  2) Slide 13 - code-listing textbox: rename loop variables "key"/"value" -> "a"/"cost"

$p = $ppt.ActivePresentation

# --- Edit 1: slide 1, date textbox ---------------------------------------
$slide1 = $p.Slides.Item(1)
$dateShape = $slide1.Shapes.Item(6)
$dateShape.TextFrame.TextRange.Text = "17. Juli 2025"

# --- Edit 2: slide 13, code textbox ---------------------------------------
$slide13 = $p.Slides.Item(13)
$codeShape = $slide13.Shapes.Item(7)
$codeTr = $codeShape.TextFrame.TextRange

# Apply the three substring replacements back-to-front so earlier offsets
# stay valid after each edit changes the overall text length.
$codeTr.Characters(499, 7).Text = " cost "
$codeTr.Characters(447, 11).Text = "a, cost "
$codeTr.Characters(426, 17).Text = "            [a "
